$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date in B1 (stored as text)
$ws.Range("B1").Value = "25/03/2023"

# Update forecast values (column B) and hour values (column C) for rows 2-20
$values = @(
    @{ Row = 2;  B = 892;   C = 14 },
    @{ Row = 3;  B = 292;   C = 14 },
    @{ Row = 4;  B = 7;     C = 14 },
    @{ Row = 5;  B = 170;   C = 14 },
    @{ Row = 6;  B = 39;    C = 14 },
    @{ Row = 7;  B = 68;    C = 14 },
    @{ Row = 8;  B = 19;    C = 14 },
    @{ Row = 9;  B = 44;    C = 14 },
    @{ Row = 10; B = 394;   C = 14 },
    @{ Row = 11; B = 218;   C = 14 },
    @{ Row = 12; B = 774;   C = 14 },
    @{ Row = 13; B = 486;   C = 14 },
    @{ Row = 14; B = 603;   C = 14 },
    @{ Row = 15; B = 167;   C = 14 },
    @{ Row = 16; B = 120;   C = 14 },
    @{ Row = 17; B = 73;    C = 14 },
    @{ Row = 18; B = 3;     C = 14 },
    @{ Row = 19; B = 9;     C = 14 },
    @{ Row = 20; B = 43;    C = 14 }
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
